$wb = $excel.ActiveWorkbook

# --- Team sheet: add Saloni's email as a mailto hyperlink in D5 ---
$team = $wb.Worksheets.Item("Team")
$team.Range("D5").Value = "ssetia@stevens.edu"
$team.Hyperlinks.Add($team.Range("D5"), "mailto:ssetia@stevens.edu")

# --- Sprint1 sheet: fill in story point estimates and update selection ---
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("E3").Value = 100
$sprint1.Range("F3").Value = 240
$sprint1.Range("E5").Value = 100
$sprint1.Range("F5").Value = 240
$null = $sprint1.Range("G10").Select()
